$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 13105.86444722955
$ws.Range("D5").Value = 13105.86444722955

$ws.Range("D9").Value = 7137.736985363312
$ws.Range("D10").Value = 7137.736985363312

$ws.Range("D14").Value = 7111.735552770167
$ws.Range("D15").Value = 7111.735552770167
